$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update revised M2 values for existing rows (open/high/low/close all
# mirror the same figure in this feed, per the source data convention) ---
$ws.Range("C122:F122").Value = 261572700000
$ws.Range("C123:F123").Value = 259192800000
$ws.Range("C124:F124").Value = 258559200000
$ws.Range("C125:F125").Value = 255274400000
$ws.Range("C131:F131").Value = 263647500000
$ws.Range("C134:F134").Value = 286126300000
$ws.Range("C135:F135").Value = 283933800000
$ws.Range("C136:F136").Value = 283623100000
$ws.Range("C137:F137").Value = 280655000000
$ws.Range("C142:F142").Value = 295604600000
$ws.Range("C143:F143").Value = 295922500000
$ws.Range("C144:F144").Value = 296506400000
$ws.Range("C146:F146").Value = 314026000000
$ws.Range("C147:F147").Value = 309824800000
$ws.Range("C148:F148").Value = 311655400000
$ws.Range("C149:F149").Value = 314511600000
$ws.Range("C150:F150").Value = 318498700000
$ws.Range("C151:F151").Value = 320618400000
$ws.Range("C152:F152").Value = 321076200000
$ws.Range("C153:F153").Value = 324218500000
$ws.Range("C154:F154").Value = 329607000000
$ws.Range("C155:F155").Value = 332065200000
$ws.Range("C156:F156").Value = 336385900000
$ws.Range("C157:F157").Value = 339280200000
$ws.Range("C158:F158").Value = 350004800000
$ws.Range("C159:F159").Value = 348823700000
$ws.Range("C160:F160").Value = 352411100000
$ws.Range("C161:F161").Value = 351238500000
$ws.Range("C162:F162").Value = 354917800000
$ws.Range("C163:F163").Value = 356539500000
$ws.Range("C164:F164").Value = 362385300000
$ws.Range("C165:F165").Value = 360221900000
$ws.Range("C170:F170").Value = 381075300000
$ws.Range("C171:F171").Value = 382602200000
$ws.Range("C172:F172").Value = 384958100000
$ws.Range("C173:F173").Value = 383090000000
$ws.Range("C182:F182").Value = 422631600000
$ws.Range("C192:F192").Value = 469280100000
$ws.Range("C194:F194").Value = 487349900000
$ws.Range("C195:F195").Value = 490302300000
$ws.Range("C196:F196").Value = 496963100000
$ws.Range("C197:F197").Value = 499199700000
$ws.Range("C206:F206").Value = 564423000000
$ws.Range("C210:F210").Value = 569711700000
$ws.Range("C212:F212").Value = 569309400000
$ws.Range("C215:F215").Value = 581768700000
$ws.Range("C218:F218").Value = 603199600000


# Append new rows 221-223, copying the style (including number format) of
# column A from the last existing row (220) so the new date cells match
# the existing "s=2" style used throughout column A.
$ws.Range("A220").Copy($ws.Range("A221"))
$ws.Range("A220").Copy($ws.Range("A222"))
$ws.Range("A220").Copy($ws.Range("A223"))

$ws.Range("A221").Value = 44986.45833333334
$ws.Range("B221").Value = "ECONOMICS:ROM2"
$ws.Range("C221:F221").Value = 613926400000
$ws.Range("G221").Value = 0

$ws.Range("A222").Value = 45017.45833333334
$ws.Range("B222").Value = "ECONOMICS:ROM2"
$ws.Range("C222:F222").Value = 618680400000
$ws.Range("G222").Value = 0

$ws.Range("A223").Value = 45047.41666666666
$ws.Range("B223").Value = "ECONOMICS:ROM2"
$ws.Range("C223:F223").Value = 624790700000
$ws.Range("G223").Value = 0

